$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(35).Insert()

$ws.Range("A35").Value = 2
$ws.Range("B35").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C35").Value = "Coquimbo"
$ws.Range("D35").Value2 = 44482
$ws.Range("E35").Value = 4
$ws.Range("F35").Value = 100112031
$ws.Range("G35").Value = "Poroto verde"
$ws.Range("H35").Value = "Magnum"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 600
$ws.Range("K35").Value = 43000
$ws.Range("L35").Value = 45000
$ws.Range("M35").Value = 44000
$ws.Range("N35").Value = "$/caja 25 kilos"
$ws.Range("O35").Value = "Provincia de Limarí"
$ws.Range("P35").Value = 1760
$ws.Range("Q35").Value = 25
$ws.Range("R35").Value = "Hortaliza"
